$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 for year 2021, copying the header-style formatting from A4 to A5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "2021年"

$ws.Range("B5").Value = 13.2
$ws.Range("C5").Value = 28.4
$ws.Range("D5").Value = -14.3
$ws.Range("E5").Value = 48.7
$ws.Range("F5").Value = -8.4
$ws.Range("H5").Value = 10.9
$ws.Range("I5").Value = 3.8
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = -0.8
$ws.Range("L5").Value = -4.1
$ws.Range("M5").Value = 10.4
$ws.Range("N5").Value = -22.8
$ws.Range("O5").Value = -6.1
$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 8.5
$ws.Range("R5").Value = 26.3
$ws.Range("S5").Value = 27.8
$ws.Range("T5").Value = 3.8
$ws.Range("U5").Value = 5.3
$ws.Range("V5").Value = -8.9
$ws.Range("W5").Value = 20
$ws.Range("X5").Value = 14.2
$ws.Range("Y5").Value = 13.1
$ws.Range("Z5").Value = 33.1
$ws.Range("AA5").Value = 14.5
$ws.Range("AB5").Value = 5.6
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 1.8
$ws.Range("AE5").Value = 11.9
$ws.Range("AF5").Value = -23.3
$ws.Range("AH5").Value = -38.6
$ws.Range("AI5").Value = -22
$ws.Range("AJ5").Value = -12.9
$ws.Range("AK5").Value = 18.7
$ws.Range("AL5").Value = -1.3
$ws.Range("AM5").Value = -3.1
$ws.Range("AN5").Value = -14
$ws.Range("AO5").Value = -15.7
$ws.Range("AP5").Value = -45.3
$ws.Range("AQ5").Value = 3.8
$ws.Range("AR5").Value = 2.8
$ws.Range("AS5").Value = 60.5
$ws.Range("AT5").Value = 114.5
$ws.Range("AU5").Value = 54.5
$ws.Range("AV5").Value = 3.8
$ws.Range("AW5").Value = 66.1
$ws.Range("AX5").Value = -13.3
$ws.Range("AY5").Value = -14.1
$ws.Range("AZ5").Value = 17
$ws.Range("BA5").Value = -2
$ws.Range("BB5").Value = 1.8
$ws.Range("BC5").Value = 12.3
$ws.Range("BD5").Value = 84.1
$ws.Range("BE5").Value = 1
$ws.Range("BF5").Value = 1.4
$ws.Range("BG5").Value = 10.8
$ws.Range("BH5").Value = -26.2
$ws.Range("BI5").Value = 6.9
$ws.Range("BJ5").Value = 13.5
$ws.Range("BK5").Value = 19.5
$ws.Range("BL5").Value = -5.5
$ws.Range("BN5").Value = -5.4
$ws.Range("BO5").Value = -4.1
$ws.Range("BP5").Value = -6.2
$ws.Range("BQ5").Value = 34.7
$ws.Range("BR5").Value = 16.6
$ws.Range("BS5").Value = -3.4
$ws.Range("BT5").Value = -1.5
$ws.Range("BU5").Value = -4.4
$ws.Range("BV5").Value = 0.4
$ws.Range("BW5").Value = 2.4
$ws.Range("BX5").Value = 28.5
$ws.Range("BY5").Value = 17.5
$ws.Range("BZ5").Value = 1.7
$ws.Range("CA5").Value = 13.8
$ws.Range("CB5").Value = 3.3
$ws.Range("CC5").Value = 14.2
$ws.Range("CD5").Value = 10.9
$ws.Range("CE5").Value = -7.6
$ws.Range("CF5").Value = 12.4
$ws.Range("CG5").Value = 10.8
$ws.Range("CH5").Value = -31
$ws.Range("CI5").Value = 9.7
$ws.Range("CJ5").Value = -13.1
$ws.Range("CK5").Value = 15.2
$ws.Range("CL5").Value = -2.4
$ws.Range("CM5").Value = -35.2
$ws.Range("CN5").Value = 22.1
$ws.Range("CO5").Value = 6
$ws.Range("CP5").Value = 17.1
$ws.Range("CQ5").Value = 11.8
$ws.Range("CR5").Value = -22.7
$ws.Range("CS5").Value = 6.5
$ws.Range("CT5").Value = 10.4
$ws.Range("CU5").Value = 12.3
$ws.Range("CV5").Value = -7.4
$ws.Range("CW5").Value = 39.2
$ws.Range("CX5").Value = 15
$ws.Range("CY5").Value = 11.8
$ws.Range("CZ5").Value = 69.8
$ws.Range("DA5").Value = 8.9
$ws.Range("DB5").Value = 5.4
$ws.Range("DC5").Value = 18.5
$ws.Range("DD5").Value = -38.1
$ws.Range("DE5").Value = -14.8
$ws.Range("DF5").Value = 17.1
$ws.Range("DG5").Value = 21.4
$ws.Range("DH5").Value = 7.6
$ws.Range("DI5").Value = 0.2
$ws.Range("DJ5").Value = 18.5
$ws.Range("DK5").Value = 27
